$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("拟合")

# Bring over the number formats from the row above (date / percent) so the
# new row matches the existing table's look instead of creating new style
# entries.
$ws.Range("A41").Copy()
$ws.Range("A42").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("F41").Copy()
$ws.Range("F42").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("J41").Copy()
$ws.Range("J42").PasteSpecial(-4122)  # xlPasteFormats

# New day of data (lockdown feature): 2022-04-10, day 41 of the outbreak.
$ws.Range("A42").Value = 44661
$ws.Range("B42").Value = 41
$ws.Range("C42").Value = 914
$ws.Range("D42").Value = 25173
$ws.Range("E42").Formula = "=D42+C42"
$ws.Range("F42").Formula = "=E42/E41-1"
$ws.Range("G42").Formula = "=SUM(C29:C42)"
$ws.Range("H42").Formula = "=SUM(D29:D42)"
$ws.Range("I42").Formula = "=SUM(E29:E42)"
$ws.Range("J42").Formula = "=I42/I41-1"

[void]$ws.Range("G42").Select()
